# Superdettagli_Light.xlsx: replace the hard-coded numeric placeholder
# (1598) in column C of the detail rows with distinct text labels, and
# re-apply the "odd row" banding style (fill + border) to all of column C
# so every row in the range looks consistent.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the already-correct formatting from C2 (fill/border used by the
# "odd" row style) onto the whole C2:C13 block, so every cell in the
# column shares the same look regardless of which row style it sits in.
$ws.Range("C2").Copy()
$ws.Range("C2:C13").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Replace the placeholder numeric values with the new per-row labels.
$ws.Range("C2").Value = "New 1"
$ws.Range("C3").Value = "New 2"
$ws.Range("C4").Value = "New 3"
$ws.Range("C5").Value = "New 4"
$ws.Range("C6").Value = "New 5"
$ws.Range("C7").Value = "New 6"
$ws.Range("C8").Value = "New 7"
$ws.Range("C9").Value = "New 8"
$ws.Range("C10").Value = "New 9"
$ws.Range("C11").Value = "New 10"
$ws.Range("C12").Value = "New 11"
$ws.Range("C13").Value = "New 12"

# Leave the selection where the author ended up after editing the range.
$ws.Range("C14").Select()
